$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextSafe($cellAddr, $value) {
    $r = $ws.Range($cellAddr)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextSafe "D2" "58.366.20"
$ws.Range("E2").Value = "  +4.13%  "
Set-TextSafe "D3" "2.454.66"
$ws.Range("E3").Value = "  +2.72%  "
$ws.Range("E4").Value = "  +0.21%  "
Set-TextSafe "D5" "159.15"
$ws.Range("E5").Value = "  +7.81%  "
Set-TextSafe "D6" "495.43"
$ws.Range("E6").Value = "  +3.50%  "
Set-TextSafe "D7" "0.616"
$ws.Range("E7").Value = "  +23.50%  "
Set-TextSafe "D8" "0.995"
$ws.Range("E8").Value = "  -0.44%  "
Set-TextSafe "D9" "2.482.76"
$ws.Range("E9").Value = "  +4.00%  "
Set-TextSafe "D10" "6.26"
$ws.Range("E10").Value = "  +14.72%  "
Set-TextSafe "D11" "0.102"
$ws.Range("E11").Value = "  +4.71%  "
$ws.Range("E12").Value = "  +4.20%  "
$ws.Range("E13").Value = "  +1.72%  "
Set-TextSafe "D14" "2.884.14"
$ws.Range("E14").Value = "  +2.86%  "
Set-TextSafe "D15" "58.252.77"
$ws.Range("E15").Value = "  +3.30%  "
Set-TextSafe "D16" "21.41"
$ws.Range("E16").Value = "  +5.25%  "
Set-TextSafe "D17" "0.0000135"
$ws.Range("E17").Value = "  +2.28%  "
Set-TextSafe "D18" "2.481.63"
$ws.Range("E18").Value = "  +4.06%  "
Set-TextSafe "D19" "4.73"
$ws.Range("E19").Value = "  +5.34%  "
Set-TextSafe "D20" "328.79"
$ws.Range("E20").Value = "  +4.59%  "
Set-TextSafe "D21" "10.15"
$ws.Range("E21").Value = "  +4.23%  "
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextSafe "D22" "0.998"
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextSafe "D23" "5.97"
$ws.Range("E23").Value = "  +5.38%  "
Set-TextSafe "D24" "58.67"
$ws.Range("E24").Value = "  +3.11%  "
Set-TextSafe "D25" "0.410"
$ws.Range("E25").Value = "  +3.96%  "
Set-TextSafe "D26" "0.164"
$ws.Range("E26").Value = "  +3.94%  "
Set-TextSafe "D27" "0.994"
$ws.Range("E27").Value = "  -0.77%  "
Set-TextSafe "D28" "2.570.37"
$ws.Range("E28").Value = "  +2.84%  "
Set-TextSafe "D29" "7.43"
$ws.Range("E29").Value = "  +1.83%  "
Set-TextSafe "D30" "0.0₃0806"
$ws.Range("E30").Value = "  +4.76%  "
Set-TextSafe "D31" "0.998"
$ws.Range("E31").Value = "  -0.14%  "
Set-TextSafe "D32" "19.02"
$ws.Range("E32").Value = "  +6.12%  "
Set-TextSafe "D33" "151.44"
$ws.Range("E33").Value = "  +1.97%  "
$ws.Range("E34").Value = "  +3.68%  "
Set-TextSafe "D35" "5.44"
$ws.Range("E35").Value = "  +9.46%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextSafe "D36" "1.17"
$ws.Range("E36").Value = "  +5.54%  "
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextSafe "D37" "3.85"
$ws.Range("E37").Value = "  +7.31%  "
Set-TextSafe "D38" "0.848"
$ws.Range("E38").Value = "  +0.08%  "
Set-TextSafe "D39" "3.64"
$ws.Range("E39").Value = "  +8.01%  "
$ws.Range("E40").Value = "  +5.15%  "
Set-TextSafe "D41" "34.39"
$ws.Range("E41").Value = "  +2.92%  "
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextSafe "D42" "284.34"
$ws.Range("E42").Value = "  +12.66%  "
$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextSafe "D43" "0.101"
$ws.Range("E43").Value = "  +7.11%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextSafe "D44" "0.992"
$ws.Range("E44").Value = "  -0.58%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextSafe "D45" "0.609"
$ws.Range("E45").Value = "  +4.19%  "
Set-TextSafe "D46" "0.0546"
$ws.Range("E46").Value = "  +0.59%  "
Set-TextSafe "D47" "0.0234"
$ws.Range("E47").Value = "  +4.71%  "
Set-TextSafe "D48" "4.76"
$ws.Range("E48").Value = "  +4.13%  "
Set-TextSafe "D49" "10.28"
$ws.Range("E49").Value = "  +0.82%  "
$ws.Range("B50").Value = "SuiNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextSafe "D50" "0.703"
$ws.Range("E50").Value = "  +14.48%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextSafe "D51" "18.14"
$ws.Range("E51").Value = "  +6.39%  "
